# Inserts two new weekly price records (for the week of 2022-02-18,
# serial 44610) into the "Fruta, Mercado Mayorista Lo Valledor de
# Santiago - Platano" dataset. The two new rows are inserted just above
# the existing row 897, pushing the rest of the table (and the two
# trailing rows) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 897 - everything below (formerly
# rows 897:957) shifts down to 899:959.
$ws.Rows("897:898").Insert()

# ---- Row 897: Platano, Sin especificar, Pinton ----
$ws.Range("A897").Value = 6
$ws.Range("B897").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C897").Value = "Metropolitana"
$ws.Range("D897").Value = 44610
$ws.Range("E897").Value = 13
$ws.Range("F897").Value = "Fruta"
$ws.Range("G897").Value = 100108
$ws.Range("H897").Value = "Tropicales y subtropicales"
$ws.Range("I897").Value = 100108006
$ws.Range("J897").Value = "Plátano"
$ws.Range("K897").Value = "Sin especificar"
$ws.Range("L897").Value = "Pintón"
$ws.Range("M897").Value = 1840
$ws.Range("N897").Value = 12000
$ws.Range("O897").Value = 13000
$ws.Range("P897").Value = 12391
$ws.Range("Q897").Value = "$/caja 20 kilos"
$ws.Range("R897").Value = "Ecuador"
$ws.Range("S897").Value = 620
$ws.Range("T897").Value = 20

# ---- Row 898: Platano, Sin especificar, Primera Pinton ----
$ws.Range("A898").Value = 6
$ws.Range("B898").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C898").Value = "Metropolitana"
$ws.Range("D898").Value = 44610
$ws.Range("E898").Value = 13
$ws.Range("F898").Value = "Fruta"
$ws.Range("G898").Value = 100108
$ws.Range("H898").Value = "Tropicales y subtropicales"
$ws.Range("I898").Value = 100108006
$ws.Range("J898").Value = "Plátano"
$ws.Range("K898").Value = "Sin especificar"
$ws.Range("L898").Value = "Primera Pintón"
$ws.Range("M898").Value = 3540
$ws.Range("N898").Value = 14000
$ws.Range("O898").Value = 16000
$ws.Range("P898").Value = 14740
$ws.Range("Q898").Value = "$/caja 20 kilos"
$ws.Range("R898").Value = "Ecuador"
$ws.Range("S898").Value = 737
$ws.Range("T898").Value = 20
